# Reverse the comma-separated list of names/emails in column G ("Recorded By")
# for every data row, except for the two specific combinations that were left
# untouched by the author ("System, dnasr281@gmail.com" and
# "System, admin@admin.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Values that must NOT be reversed even though they contain a comma.
$exempt = @("System, dnasr281@gmail.com", "System, admin@admin.com")

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Length -le 1) {
        continue
    }

    $isExempt = $false
    foreach ($ex in $exempt) {
        if ($text -eq $ex) {
            $isExempt = $true
        }
    }

    if ($isExempt) {
        continue
    }

    $reversedParts = $parts[($parts.Length - 1)..0]
    $newText = [string]::Join(", ", $reversedParts)

    $cell.Value = $newText
}
